# Applies the "end of loop" progress update to the contractor bill sheet:
# refresh the executed-quantity figures (column C) and the corresponding
# "Upto date Amount" text values (column G), as well as the Grand Total /
# Net Payable totals that roll those amounts up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Qty executed upto date (plain numeric cells) ---
$ws.Range("C8").Value  = 16
$ws.Range("C9").Value  = 66
$ws.Range("C10").Value = 79
$ws.Range("C11").Value = 58
$ws.Range("C12").Value = 32
$ws.Range("C13").Value = 13
$ws.Range("C14").Value = 70
$ws.Range("C15").Value = 76
$ws.Range("C16").Value = 6
$ws.Range("C17").Value = 68

# Helper: write a value into a cell as TEXT (matching the workbook's existing
# convention of storing "Upto date Amount" figures as formatted text strings
# like "16896.00") without leaving behind any new/changed cell style.
function Set-TextValue {
    param($cell, [string]$text)

    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $originalStyle
}

# --- Upto date Amount (text values formatted as "N.00") ---
Set-TextValue $ws.Range("G9")  "16896.00"
Set-TextValue $ws.Range("G10") "37288.00"
Set-TextValue $ws.Range("G11") "38396.00"
Set-TextValue $ws.Range("G13") "1768.00"
Set-TextValue $ws.Range("G14") "1610.00"

# --- Grand Total rows ---
Set-TextValue $ws.Range("G19") "95958.00"
Set-TextValue $ws.Range("H19") "95958.00"
Set-TextValue $ws.Range("G21") "95958.00"
Set-TextValue $ws.Range("H21") "95958.00"
